# Auto-generated edits applying scheduled-runner price/profit recalculations
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2359306.2
$ws.Range("I17").Value = 163
$ws.Range("J17").Value = 2394517.2
$ws.Range("K17").Value = 489
$ws.Range("L17").Value = 7183551.600000001
$ws.Range("M17").Value = -321
$ws.Range("N17").Value = -7183887.600000001

$ws.Range("H19").Value = 450.45456
$ws.Range("I19").Value = 487.5
$ws.Range("J19").Value = 429.2857
$ws.Range("K19").Value = 487.5
$ws.Range("L19").Value = 429.2857
$ws.Range("M19").Value = -312.5
$ws.Range("N19").Value = -779.2857

$ws.Range("H38").Value = 1595.7059
$ws.Range("I38").Value = 179.75
$ws.Range("J38").Value = 2854.3333
$ws.Range("K38").Value = 539.25
$ws.Range("L38").Value = 8562.999899999999
$ws.Range("M38").Value = -167.25
$ws.Range("N38").Value = -9306.999899999999

$ws.Range("H39").Value = 574.7
$ws.Range("I39").Value = 133.33333
$ws.Range("J39").Value = 763.8570999999999
$ws.Range("K39").Value = 399.99999
$ws.Range("L39").Value = 2291.5713
$ws.Range("M39").Value = -103.99999
$ws.Range("N39").Value = -2883.5713

$ws.Range("H51").Value = 6857
$ws.Range("I51").Value = 2999.5
$ws.Range("J51").Value = 8400
$ws.Range("K51").Value = 2999.5
$ws.Range("L51").Value = 8400
$ws.Range("M51").Value = -2515.5
$ws.Range("N51").Value = -9368

$ws.Range("H74").Value = 4349.9165
$ws.Range("I74").Value = 4466.6665
$ws.Range("J74").Value = 3999.6667
$ws.Range("K74").Value = 4466.6665
$ws.Range("L74").Value = 3999.6667
$ws.Range("M74").Value = -3530.6665
$ws.Range("N74").Value = -5871.6667

$ws.Range("H77").Value = 4349.9165
$ws.Range("I77").Value = 4466.6665
$ws.Range("J77").Value = 3999.6667
$ws.Range("K77").Value = 22333.3325
$ws.Range("L77").Value = 19998.3335
$ws.Range("M77").Value = -17653.3325
$ws.Range("N77").Value = -29358.3335

$ws.Range("H128").Value = 29610.824
$ws.Range("J128").Value = 29610.824
$ws.Range("L128").Value = 29610.824
$ws.Range("N128").Value = -39570.824

$ws.Range("H132").Value = 1363.5454
$ws.Range("I132").Value = 1404.6666
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 4213.9998
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -1683.9998
$ws.Range("N132").Value = -6560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 11241.1
$ws.Range("I45").Value = 12334.556
$ws.Range("J45").Value = 1400
$ws.Range("K45").Value = 12334.556
$ws.Range("L45").Value = 1400
$ws.Range("M45").Value = -11957.556
$ws.Range("N45").Value = -2154

$ws.Range("H74").Value = 1555.7084
$ws.Range("I74").Value = 1556.5238
$ws.Range("K74").Value = 1556.5238
$ws.Range("M74").Value = -682.5237999999999

$ws.Range("H77").Value = 1555.7084
$ws.Range("I77").Value = 1556.5238
$ws.Range("K77").Value = 7782.619
$ws.Range("M77").Value = -3414.619

$ws.Range("H122").Value = 1285414.5
$ws.Range("I122").Value = 1511552.4
$ws.Range("J122").Value = 3966.3333
$ws.Range("K122").Value = 4534657.199999999
$ws.Range("L122").Value = 11898.9999
$ws.Range("M122").Value = -4532207.199999999
$ws.Range("N122").Value = -16798.9999

$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 50000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("M123").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6802.4517
$ws.Range("I31").Value = 1887.85
$ws.Range("J31").Value = 15738.091
$ws.Range("K31").Value = 1887.85
$ws.Range("L31").Value = 15738.091
$ws.Range("M31").Value = -1592.85
$ws.Range("N31").Value = -16328.091

$ws.Range("H34").Value = 6802.4517
$ws.Range("I34").Value = 1887.85
$ws.Range("J34").Value = 15738.091
$ws.Range("K34").Value = 1887.85
$ws.Range("L34").Value = 15738.091
$ws.Range("M34").Value = -1685.85
$ws.Range("N34").Value = -16142.091

$ws.Range("H132").Value = 3377.4285
$ws.Range("I132").Value = 3077.2942
$ws.Range("J132").Value = 4653
$ws.Range("K132").Value = 9231.882599999999
$ws.Range("L132").Value = 13959
$ws.Range("M132").Value = -6701.882599999999
$ws.Range("N132").Value = -19019

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 20000152
$ws.Range("I23").Value = 100000000
$ws.Range("J23").Value = 191
$ws.Range("K23").Value = 300000000
$ws.Range("L23").Value = 573
$ws.Range("M23").Value = -299999765
$ws.Range("N23").Value = -1043

$ws.Range("H113").Value = 2500581
$ws.Range("I113").Value = 4167154
$ws.Range("J113").Value = 834007.9399999999
$ws.Range("K113").Value = 12501462
$ws.Range("L113").Value = 2502023.82
$ws.Range("M113").Value = -12499292
$ws.Range("N113").Value = -2506363.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 70004
$ws.Range("J4").Value = 70004
$ws.Range("L4").Value = 70004
$ws.Range("N4").Value = -70228

$ws.Range("H107").Value = 1989
$ws.Range("I107").Value = 741.2857
$ws.Range("J107").Value = 2716.8333
$ws.Range("K107").Value = 741.2857
$ws.Range("L107").Value = 2716.8333
$ws.Range("M107").Value = 1178.7143
$ws.Range("N107").Value = -6556.8333

$ws.Range("H122").Value = 3589663.2
$ws.Range("I122").Value = 3087948.8
$ws.Range("J122").Value = 4547482
$ws.Range("K122").Value = 9263846.399999999
$ws.Range("L122").Value = 13642446
$ws.Range("M122").Value = -9261396.399999999
$ws.Range("N122").Value = -13647346

$ws.Range("H132").Value = 3868
$ws.Range("I132").Value = 3733.1765
$ws.Range("J132").Value = 4002.8235
$ws.Range("K132").Value = 11199.5295
$ws.Range("L132").Value = 12008.4705
$ws.Range("M132").Value = -8669.529500000001
$ws.Range("N132").Value = -17068.4705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 3995
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 5990
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 5990
$ws.Range("M34").Value = -1828
$ws.Range("N34").Value = -6334

$ws.Range("H46").Value = 16667639
$ws.Range("I46").Value = 37037816
$ws.Range("J46").Value = 1130.8182
$ws.Range("K46").Value = 37037816
$ws.Range("L46").Value = 1130.8182
$ws.Range("M46").Value = -37037628
$ws.Range("N46").Value = -1506.8182

$ws.Range("H122").Value = 3879882.8
$ws.Range("I122").Value = 4466752.5
$ws.Range("J122").Value = 2001900
$ws.Range("K122").Value = 13400257.5
$ws.Range("L122").Value = 6005700
$ws.Range("M122").Value = -13397807.5
$ws.Range("N122").Value = -6010600

$ws.Range("H132").Value = 14450918
$ws.Range("I132").Value = 17340058
$ws.Range("J132").Value = 5218.6
$ws.Range("K132").Value = 52020174
$ws.Range("L132").Value = 15655.8
$ws.Range("M132").Value = -52017644
$ws.Range("N132").Value = -20715.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1806.7693
$ws.Range("I122").Value = 1728.8
$ws.Range("J122").Value = 2066.6667
$ws.Range("K122").Value = 5186.4
$ws.Range("L122").Value = 6200.000100000001
$ws.Range("M122").Value = -2736.4
$ws.Range("N122").Value = -11100.0001

$ws.Range("H132").Value = 1551.3636
$ws.Range("I132").Value = 1062.9445
$ws.Range("J132").Value = 3749.25
$ws.Range("K132").Value = 3188.8335
$ws.Range("L132").Value = 11247.75
$ws.Range("M132").Value = -658.8335000000002
$ws.Range("N132").Value = -16307.75
